# Situation Report - periodic refresh:
#   1. Bump the report date "02 September, 2021" -> "03 September, 2021".
#   2. Regenerate the keys (bookmarks) that back the two table
#      cross-references ("tab:OverviewTable" / "tab:StateLevelTable") so
#      they get freshly minted identities, per the new table-keys
#      generation function.

$d = $word.ActiveDocument

# --- 1. Date bump -----------------------------------------------------
# Only the leading "02" token (the day-of-month) in the Date-styled
# paragraph should change; everything else in that paragraph (and the
# "2021"/"0.02" etc. substrings elsewhere in the document) must stay
# untouched, so we scope the edit to an exact 2-character Range rather
# than doing a document-wide Find/Replace.
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Date") {
        $dayRange = $d.Range($p.Range.Start, $p.Range.Start + 2)
        if ($dayRange.Text -eq "02") {
            $dayRange.Delete()
            $d.Range($p.Range.Start, $p.Range.Start).InsertBefore("03")
        }
    }
}

# --- 2. Refresh the table-caption bookmark keys ------------------------
function Refresh-TableKey([string]$name) {
    $bm = $d.Bookmarks.Item($name)
    $rng = $bm.Range
    $bm.Delete()
    $d.Bookmarks.Add($name, $rng)
}

Refresh-TableKey "tab:OverviewTable"
Refresh-TableKey "tab:StateLevelTable"
